$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume 1h (E) columns,
# plus a few re-ordered rows where Coin name (B) and Link (C) changed.
# Values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the original
# inlineStr cells, e.g. '8.00' must not become the number 8), then
# the style is reset to Normal so no numeric formatting sticks.

$ws.Range('D2').Value = '67.654.36'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '3.503.08'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = "'606.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.47%  '
$ws.Range('D6').Value = "'150.11"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('D7').Value = '3.501.50'
$ws.Range('E7').Value = '  -3.28%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('D11').Value = "'7.56"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('D12').Value = "'0.429"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.49%  '
$ws.Range('E13').Value = '  -5.17%  '
$ws.Range('D14').Value = "'32.09"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').Value = '4.087.54'
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '67.751.03'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.481.38'
$ws.Range('E17').Value = '  -3.58%  '
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = "'6.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.88%  '
$ws.Range('D20').Value = "'15.53"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.70%  '
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').Value = "'448.84"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').Value = "'0.625"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = "'78.85"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').Value = '3.645.67'
$ws.Range('E25').Value = '  -3.16%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = "'0.0000123"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -9.13%  '
$ws.Range('D28').Value = "'8.73"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').Value = "'9.93"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.76%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = "'1.66"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.51"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.36%  '
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'6.21"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.68%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = "'25.58"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.38%  '
$ws.Range('D36').Value = '3.494.82'
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('E37').Value = '  -6.22%  '
$ws.Range('D38').Value = "'8.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.92%  '
$ws.Range('E40').Value = '  -5.28%  '
$ws.Range('D41').Value = "'178.32"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').Value = "'0.997"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = "'0.0904"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').Value = "'31.01"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('D46').Value = "'0.899"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').Value = "'46.88"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('E48').Value = '  -3.93%  '
$ws.Range('E49').Value = '  -2.10%  '
$ws.Range('D50').Value = "'2.50"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -10.49%  '
$ws.Range('D51').Value = "'0.999"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.60%  '
